# [EI-979] Survey.xlsx data dictionary header update.
# Renames the Sheet1 header cells (row 1) from "Then_Question"/"Else_Question"
# to "Then_Goto"/"Else_Goto" (columns I and J), and updates the saved
# selection/active cell on Sheet1 to I1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the header text for the goto columns (I1 / J1).
$ws.Range("I1").Value = "Then_Goto"
$ws.Range("J1").Value = "Else_Goto"

# Update the saved selection/active cell on Sheet1 to I1.
$ws.Range("I1").Select()

$wb.Save()
